$d = $word.ActiveDocument

# 1. Remove the existing hidden "_GoBack" bookmark that currently sits
#    between "hope" and " I think)" in the second paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Append two new paragraphs at the end of the document body:
#      "Iniquity or lawlessness"
#      "Ordinances or statutes" + " or truths" (as two separate runs),
#    followed immediately by the "_GoBack" bookmark (re-homed to the new
#    end of the document).
#    InsertXML is used so the run boundaries are preserved exactly as
#    authored (a plain Range.Text/InsertAfter merges same-format runs
#    together) and so the bookmark lands precisely at end-of-paragraph
#    without relying on a collapsed Range positioned on the paragraph
#    mark (Bookmarks.Add mis-resolves a Range collapsed exactly there).
$endRange = $d.Content
$endRange.Collapse(0)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Iniquity or lawlessness</w:t></w:r></w:p><w:p><w:r><w:t>Ordinances or statutes</w:t></w:r><w:r><w:t xml:space="preserve"> or truths</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$endRange.InsertXML($xml)
